$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.700.14"

$ws.Range("D3").Value = "'3.468.61"
$ws.Range("E3").Value = "  +1.23%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'414.80"
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("D6").Value = "'130.82"
$ws.Range("E6").Value = "  +1.71%  "

$ws.Range("E7").Value = "  -0.76%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -0.91%  "

$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("D11").Value = "'42.89"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").Value = "'9.77"
$ws.Range("E12").Value = "  +7.12%  "

$ws.Range("D13").Value = "'0.0000219"
$ws.Range("E13").Value = "  +5.00%  "

$ws.Range("D14").Value = "'4.014.95"
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("E15").Value = "  -0.23%  "

$ws.Range("D16").Value = "'20.53"
$ws.Range("E16").Value = "  -4.22%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.489.87"
$ws.Range("E17").Value = "  +1.75%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.84"
$ws.Range("E18").Value = "  +2.98%  "

$ws.Range("D20").Value = "'62.659.96"
$ws.Range("E20").Value = "  +0.98%  "

$ws.Range("D21").Value = "'471.16"
$ws.Range("E21").Value = "  +5.60%  "

$ws.Range("D22").Value = "'90.93"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("E23").Value = "  +2.49%  "

$ws.Range("E24").Value = "  +3.30%  "

$ws.Range("D25").Value = "'10.56"
$ws.Range("E25").Value = "  +21.07%  "

$ws.Range("E26").Value = "  +2.28%  "

$ws.Range("D27").Value = "'33.42"
$ws.Range("E27").Value = "  +1.91%  "

$ws.Range("D28").Value = "'4.81"
$ws.Range("E28").Value = "  +0.84%  "

$ws.Range("E29").Value = "  -2.46%  "

$ws.Range("D30").Value = "'12.06"
$ws.Range("E30").Value = "  +0.38%  "

$ws.Range("D31").Value = "'2.65"
$ws.Range("E31").Value = "  -2.86%  "

$ws.Range("E32").Value = "  -2.58%  "

$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").Value = "'41.15"
$ws.Range("E34").Value = "  -4.23%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Value = "'58.92"
$ws.Range("E36").Value = "  +9.78%  "

$ws.Range("D37").Value = "'0.0490"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("E39").Value = "  +3.93%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'3.37"
$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.323"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("E43").Value = "  +6.43%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'4.39"
$ws.Range("E44").Value = "  +3.27%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'144.62"
$ws.Range("E45").Value = "  +1.91%  "

$ws.Range("E46").Value = "  +4.23%  "

$ws.Range("E47").Value = "  +11.18%  "

$ws.Range("D48").Value = "'0.0₃0564"
$ws.Range("E48").Value = "  +38.53%  "

$ws.Range("D49").Value = "'16.41"
$ws.Range("E49").Value = "  -1.31%  "

$ws.Range("D50").Value = "'22.38"
$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("E51").Value = "  -1.23%  "

